# "risultati L1 500 immagini"
# - Fill in the previously-empty "500 scansioni (L1)" results row (row 17).
# - Remove the now-unused "700 scansioni (L1)/(LCS)" rows (20:21) content and
#   the two trailing blank rows (22:23), shrinking the used range from
#   A2:S23 down to A2:S21.
# - Update the view's active cell / scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill row 17 ("500 scansioni (L1)") with the measured results ---------
$ws.Range("D17").Value = 18146
$ws.Range("F17").Value = 3957
$ws.Range("H17").Value = 289.33
$ws.Range("I17").Value = 4482.02
$ws.Range("K17").Value = 85545.11
$ws.Range("L17").Value = 2854
$ws.Range("M17").Value = 7354
$ws.Range("P17").Value = 90.66
$ws.Range("Q17").Value = 80.89

# --- Clear out the old "700 scansioni" rows (20:21), reusing the same ----
# --- blank formatting already used by rows 18:19 --------------------------
$ws.Range("C18:Q19").Copy()
$ws.Range("C20:Q21").PasteSpecial(-4122)
$ws.Range("A20:Q21").ClearContents()
$ws.Range("A20:B21").Clear()

# --- Drop the trailing blank rows entirely --------------------------------
$ws.Rows("22:23").Delete()

# --- Update the view: scroll right so column D is leftmost, select N17 ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("N17").Select()
